# Inserts a new weekly price record as row 98 in the "Zapallo italiano"
# price sheet, pushing the existing rows 98:187 down to 99:188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 98 (shifts rows 98-187 -> 99-188)
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new record's data
$ws.Cells.Item(98, 1).Value = 4
$ws.Cells.Item(98, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(98, 3).Value = "Los Lagos"
$ws.Cells.Item(98, 4).Value = 44589
$ws.Cells.Item(98, 5).Value = 10
$ws.Cells.Item(98, 6).Value = 100112032
$ws.Cells.Item(98, 7).Value = "Zapallo italiano"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 200
$ws.Cells.Item(98, 11).Value = 16000
$ws.Cells.Item(98, 12).Value = 16000
$ws.Cells.Item(98, 13).Value = 16000
$ws.Cells.Item(98, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 320
$ws.Cells.Item(98, 17).Value = 50
$ws.Cells.Item(98, 18).Value = "Hortaliza"
